$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column G, shifting odp/smg/Reference Unit/Reference to the right
$ws.Columns("G:G").Insert()

# New header for inserted column
$ws.Range("G1").Value = "stored_carbon"

# New values for inserted column (all zero, no special number format)
$ws.Range("G2").Value = 0
$ws.Range("G3").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("G5").Value = 0

# Update the active selection to match the post-edit state
$ws.Range("G2:G5").Select()
